$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8443170964660937
$ws.Range("B3").Value = 0.8368
$ws.Range("B4").Value = 0.8955479452054794
$ws.Range("B5").Value = 0.8651778329197684
